# New API Query - 2023 Included
# API query to UN performed 11/26/2023. Query modified to include 2023 data.
#
# - short-url value (B2) refreshed to the new query token
# - "oip" (U2) and "hst" (V2) now report "-" instead of "null"/0,
#   and V2 picks up the same (left-aligned) style U2 already uses.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "3A0m7u"
$ws.Range("U2").Value = "-"
$ws.Range("V2").Value = "-"
$ws.Range("V2").HorizontalAlignment = $ws.Range("U2").HorizontalAlignment
